$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 5.5
$ws.Range("K2").Value = 1.69
$ws.Range("M2").Value = 1.2
$ws.Range("N2").Value = 4.33
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 1.22
$ws.Range("S2").Value = 10
$ws.Range("T2").Value = 1.06
$ws.Range("U2").Value = 1.88
$ws.Range("V2").Value = 1.93
$ws.Range("AA2").Value = 12
$ws.Range("AB2").Value = 17
$ws.Range("AE2").Value = 4.33
$ws.Range("AI2").Value = 8.5
$ws.Range("AM2").Value = 67

$ws.Range("G3").Value = 1.48
$ws.Range("H3").Value = 3.6
$ws.Range("I3").Value = 9.5
$ws.Range("J3").Value = 2.1
$ws.Range("K3").Value = 1.95
$ws.Range("L3").Value = 9.5
$ws.Range("M3").Value = 1.14
$ws.Range("N3").Value = 5.5
$ws.Range("O3").Value = 1.57
$ws.Range("P3").Value = 2.25
$ws.Range("Q3").Value = 2.7
$ws.Range("R3").Value = 1.44
$ws.Range("S3").Value = 6
$ws.Range("T3").Value = 1.13
$ws.Range("U3").Value = 1.62
$ws.Range("V3").Value = 2.2
$ws.Range("W3").Value = 3.25
$ws.Range("X3").Value = 1.33
$ws.Range("Z3").Value = 5
$ws.Range("AB3").Value = 9
$ws.Range("AC3").Value = 19
$ws.Range("AE3").Value = 5.5
$ws.Range("AF3").Value = 8
$ws.Range("AG3").Value = 34
$ws.Range("AH3").Value = 151
$ws.Range("AI3").Value = 13
$ws.Range("AJ3").Value = 41
$ws.Range("AK3").Value = 34
$ws.Range("AL3").Value = 126
$ws.Range("AM3").Value = 101
$ws.Range("AN3").Value = 126
$ws.Range("AP3").Value = 2.1
$ws.Range("AQ3").Value = 1.78
